# Edit script: add "Cart" node to the rentCar component tree diagram
# (sibling of "Cars" under "Main", plus its own child count-node),
# re-laying out the right-hand branch of the tree to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reset the diagram area (unmerge + clear old content/formatting) ---
$ws.Range("I1:W8").UnMerge()
$ws.Range("I1:Z9").Clear()

# --- Column widths: extend formatted columns from A:AN to A:AQ ---
# (Columns 1-40 already carry the 3.7109375 custom width from the template;
#  widen the newly-used columns 41-43 to match so the whole tree area is
#  uniformly sized, same as the author's manual column selection.)
$ws.Columns("AO:AQ").ColumnWidth = 3.7109375

# --- Set text values (tree node labels) ---
$ws.Range("Q1").Value = "APP"
$ws.Range("Q3").Value = "Header"
$ws.Range("Y3").Value = "Main"
$ws.Range("O6").Value = "Nav"
$ws.Range("R6").Value = "CarLogo"
$ws.Range("V6").Value = "Cart"
$ws.Range("Y6").Value = "Cars"
$ws.Range("I9").Value = "Logo"
$ws.Range("L9").Value = "DarkMode"
$ws.Range("Q9").Value = "MenuToggle"
$ws.Range("V9").Value = "Cart"

# --- Merge the tree-node label cells ---
$ws.Range("L9:O9").Merge()
$ws.Range("I9:J9").Merge()
$ws.Range("Q3:R3").Merge()
$ws.Range("Y3:Z3").Merge()
$ws.Range("Q1:Z1").Merge()
$ws.Range("O6:P6").Merge()
$ws.Range("R6:S6").Merge()
$ws.Range("Q9:T9").Merge()
$ws.Range("Y6:Z6").Merge()
$ws.Range("V6:W6").Merge()
$ws.Range("V9:W9").Merge()

# --- Apply borders/alignment per group of cells sharing the same style ---
foreach ($addr in @("Q4", "Y4", "Q6", "S7", "T7", "U7", "V7", "Y7", "Z7", "U9")) {
  $rng = $ws.Range($addr)
  $rng.HorizontalAlignment = -4108
}

foreach ($addr in @("S4", "T4", "U4", "V4", "K7", "L7", "M7", "N7", "R7")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(9).LineStyle = 1
}

foreach ($addr in @("O7", "Q7")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(9).LineStyle = 1
  $rng.HorizontalAlignment = -4108
}

foreach ($addr in @("Q5")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(8).LineStyle = 1
}

foreach ($addr in @("M9", "N9")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(8).LineStyle = 1
  $rng.Borders.Item(9).LineStyle = 1
  $rng.HorizontalAlignment = -4108
}

foreach ($addr in @("Q2")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(10).LineStyle = 1
}

foreach ($addr in @("M8")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(10).LineStyle = 1
  $rng.Borders.Item(9).LineStyle = 1
}

foreach ($addr in @("R5")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(10).LineStyle = 1
  $rng.Borders.Item(8).LineStyle = 1
  $rng.Borders.Item(9).LineStyle = 1
}

foreach ($addr in @("J9", "O9")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(10).LineStyle = 1
  $rng.Borders.Item(8).LineStyle = 1
  $rng.Borders.Item(9).LineStyle = 1
  $rng.HorizontalAlignment = -4108
}

foreach ($addr in @("Z2", "Z5")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(7).LineStyle = 1
}

foreach ($addr in @("Z4")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(7).LineStyle = 1
  $rng.HorizontalAlignment = -4108
}

foreach ($addr in @("W5", "S8", "W8")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(7).LineStyle = 1
  $rng.Borders.Item(9).LineStyle = 1
}

foreach ($addr in @("W7")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(7).LineStyle = 1
  $rng.Borders.Item(8).LineStyle = 1
  $rng.HorizontalAlignment = -4108
}

foreach ($addr in @("P5", "J8")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(7).LineStyle = 1
  $rng.Borders.Item(8).LineStyle = 1
  $rng.Borders.Item(9).LineStyle = 1
}

foreach ($addr in @("R4", "P7", "I9", "L9")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(7).LineStyle = 1
  $rng.Borders.Item(8).LineStyle = 1
  $rng.Borders.Item(9).LineStyle = 1
  $rng.HorizontalAlignment = -4108
}

foreach ($addr in @("Q1", "R1", "S1", "T1", "U1", "V1", "W1", "X1", "Y1", "Z1", "Q3", "R3", "Y3", "Z3", "O6", "P6", "R6", "S6", "V6", "W6", "Y6", "Z6", "Q9", "R9", "S9", "T9", "V9", "W9")) {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(7).LineStyle = 1
  $rng.Borders.Item(10).LineStyle = 1
  $rng.Borders.Item(8).LineStyle = 1
  $rng.Borders.Item(9).LineStyle = 1
  $rng.HorizontalAlignment = -4108
}

# --- Selection matches the author's final cursor position ---
$ws.Range("V9:W9").Select()
